$d = $word.ActiveDocument

# --- 1. Remove the "Meta description: ..." paragraph that follows the title ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- 2. Insert a new bold paragraph right before the last paragraph
#        (the one holding the image-prompt / meta-description text) ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

# The freshly inserted paragraph is now the last-but-one; re-fetch it by
# index since the earlier reference does not track the new content.
$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Florageddon! Slot Free - Review of Gameplay &amp; Winning Potential</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml) | Out-Null

# --- 3. Replace the text of the final paragraph (the old image-prompt
#        text) with the new meta-description text, keeping its italic
#        run formatting intact ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = "Read a review of Florageddon! slot game, with features, gameplay, winning potential and device availability. Play free now."
